# Insert a new weekly data row before row 311, shifting existing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(311).Insert()

$ws.Cells.Item(311, 1).Value = 11
$ws.Cells.Item(311, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(311, 3).Value = "Bíobío"
$ws.Cells.Item(311, 4).Value = 45223
$ws.Cells.Item(311, 5).Value = 8
$ws.Cells.Item(311, 6).Value = 100112003
$ws.Cells.Item(311, 7).Value = "Ajo"
$ws.Cells.Item(311, 8).Value = "Chino"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 180
$ws.Cells.Item(311, 11).Value = 19000
$ws.Cells.Item(311, 12).Value = 20000
$ws.Cells.Item(311, 13).Value = 19556
$ws.Cells.Item(311, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(311, 15).Value = "China"
$ws.Cells.Item(311, 16).Value = 1956
$ws.Cells.Item(311, 17).Value = 10
$ws.Cells.Item(311, 18).Value = "Hortaliza"
